# Refined metadata to be additional tab
#
# 1) Update the "panel_query_time" (F column) timestamps on the "data"
#    sheet to reflect the re-run query time.
# 2) Add a new "metadata" worksheet (after "data") summarising the panel
#    query itself.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- 1. refresh the time_taken column on the data sheet ------------------
$ws1.Range("F2").Value  = "2021-10-05 14:33:38.360328"
$ws1.Range("F3").Value  = "2021-10-05 14:33:38.360337"
$ws1.Range("F4").Value  = "2021-10-05 14:33:38.360340"
$ws1.Range("F5").Value  = "2021-10-05 14:33:38.360343"
$ws1.Range("F6").Value  = "2021-10-05 14:33:38.360346"
$ws1.Range("F7").Value  = "2021-10-05 14:33:38.360349"
$ws1.Range("F8").Value  = "2021-10-05 14:33:38.360352"
$ws1.Range("F9").Value  = "2021-10-05 14:33:38.360355"
$ws1.Range("F10").Value = "2021-10-05 14:33:38.360358"
$ws1.Range("F11").Value = "2021-10-05 14:33:38.360361"
$ws1.Range("F12").Value = "2021-10-05 14:33:38.360364"
$ws1.Range("F13").Value = "2021-10-05 14:33:38.360367"
$ws1.Range("F14").Value = "2021-10-05 14:33:38.360370"
$ws1.Range("F15").Value = "2021-10-05 14:33:38.360372"
$ws1.Range("F16").Value = "2021-10-05 14:33:38.360375"
$ws1.Range("F17").Value = "2021-10-05 14:33:38.360378"
$ws1.Range("F18").Value = "2021-10-05 14:33:38.360381"
$ws1.Range("F19").Value = "2021-10-05 14:33:38.360384"
$ws1.Range("F20").Value = "2021-10-05 14:33:38.360387"
$ws1.Range("F21").Value = "2021-10-05 14:33:38.360390"
$ws1.Range("F22").Value = "2021-10-05 14:33:38.360392"
$ws1.Range("F23").Value = "2021-10-05 14:33:38.360395"
$ws1.Range("F24").Value = "2021-10-05 14:33:38.360398"
$ws1.Range("F25").Value = "2021-10-05 14:33:38.360400"
$ws1.Range("F26").Value = "2021-10-05 14:33:38.360403"
$ws1.Range("F27").Value = "2021-10-05 14:33:38.360406"
$ws1.Range("F28").Value = "2021-10-05 14:33:38.360409"
$ws1.Range("F29").Value = "2021-10-05 14:33:38.360412"

# --- 2. add the "metadata" sheet, placed after "data" --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Corneal Dystrophy"
$ws2.Range("C2").Value = 91
# "data_version" is stored as text ("1.5"), not a number - a leading
# apostrophe forces Excel to keep it as text rather than coercing it to
# a numeric cell; reset back to the default style afterwards so the
# quote-prefix formatting flag doesn't linger on the cell.
$ws2.Range("D2").Value = "'1.5"
$ws2.Range("D2").Style = "Normal"
$ws2.Range("E2").Value = "2021-08-28T08:14:48.139321Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:38.356583"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/91/?format=json"

# match the header / row-index styling used on the "data" sheet (bold,
# centered, bordered cells - style index 1 there)
$ws1.Range("B1:F1").Copy() | Out-Null
$ws2.Range("B1:G1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null

$ws2.Range("A1").Select()

# keep "data" as the active sheet/tab, as before the edit
$ws1.Activate()
